# The workbook holds 16 worksheets ("23" down to "8"), each with a
# statsmodels OLS "Regression Results" text dump in cell B2. Every one of
# those dumps was re-generated a day later ("Wed, 01 Jan 2020" / "23:18:56"
# -> "Thu, 02 Jan 2020" / "20:48:49") while all the numeric results stayed
# identical - so we just patch the Date/Time banner line in place on every
# sheet, preserving the fixed-width column alignment.

$wb = $excel.ActiveWorkbook

$oldDate = "Wed, 01 Jan 2020"
$newDate = "Thu, 02 Jan 2020"
$oldTime = "23:18:56"
$newTime = "20:48:49"

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $cell = $ws.Range("B2")
    $text = $cell.Value2

    if ($text -ne $null -and $text.IndexOf($oldDate) -ge 0) {
        $updated = $text.Replace($oldDate, $newDate).Replace($oldTime, $newTime)
        $cell.Value = $updated
    }
}
